# Appends a fresh scrape run (2025-09-18) duplicating the prior run's 42 rows
# onto the "erros_carrefour_2025-09" worksheet, growing it from A1:D177 to A1:D219.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Páprica Defumada Kitano Reserva 36g", "https://mercado.carrefour.com.br/farofa-de-mandioca-tradicional-yoki-400g-6582613/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/massa-para-pastel-discao-massa-leve-500g-841757/p"),
    @("Pimentão Block Vermelho Trebeshi 150 g", "https://mercado.carrefour.com.br/pimentao-block-vermelho-trebeshi-150-g-5738458/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/chocolate-ao-leite-com-amendoim-shot-165g-5790859/p"),
    @("Couve Flor Cledson 300 g", "https://mercado.carrefour.com.br/couve-flor-cledson-300-g-9560297/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/patinho-fracionado-a-vacuo-500g-18325/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/paleta-bovina-a-vacuo-500gnao-reativarcodigo-de-compra-20745/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/costela-minga-bovina-cong-aprox-2kg-224006/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/presunto-cozido-sem-capa-fatiado-aurora-aproximadamente-200-g-49450/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/mortadela-defumada-sadia-280g-5447045/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/queijo-minas-frescal-aurora-450-g-6264693/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/queijo-coalho-bom-leite-500-g-4305054/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/leite-uht-integral-piratininga-1-l-665017/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/iogurte-natural-tradicional-batavo-170g-5150439/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/manteiga-com-sal-aviacao-200-g-10010/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/creme-de-leite-ultrapasteurizado-itambe-200-g-5988921/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/requeijao-cremoso-aviacao-tradicional-220-g-10000/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/acucar-cristal-carrefour-1kg-5147300/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/mel-com-cacau-e-avela-400-g-4510146/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/geleia-de-goiaba-selecoes-c-pedacos-260-g-1280815/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/suco-de-uva-integral-maric-1-l-3538256/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/vinho-tinto-fino-seco-cabernet-sauvignon-pergola-750ml-1521709/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/whisky-red-label-johnnie-walker-1-litro-2719/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/refrigerante-coca-cola-sabor-cola-1-5-l-11087/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/cafe-torrado-e-moido-extraforte-melitta-500g-271203/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/farinha-de-trigo-dona-benta-tradicional-1kg-196416/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/azeite-extravirgem-portugues-oliveira-da-serra-500-ml-4526108/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/oleo-de-soja-soya-900ml-482616/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/margarina-qualy-com-sal-250g-4815618/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/arroz-branco-longofino-tipo-1-tio-joao-1kg-115658/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/feijao-preto-tipo-1-kicaldo-1kg-466510/p"),
    @("Arroz Branco Carrefour Classic Olimpíadas 1Kg", "https://mercado.carrefour.com.br/arroz-branco-carrefour-classic-olimpiadas-1kg-3433455/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/busca/pao%20frances"),
    @("Não encontrado", "https://mercado.carrefour.com.br/busca/biscoito%20doce"),
    @("Não encontrado", "https://mercado.carrefour.com.br/frango-inteiro-temperado-seara-assa-facil-aprox-19kg-170739/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/busca/cafe%20moido"),
    @("Não encontrado", "https://mercado.carrefour.com.br/busca/costela?page=1"),
    @("Não encontrado", "https://mercado.carrefour.com.br/costela-de-cordeiro-a-vacuo-28738/p"),
    @("Não encontrado", "https://mercado.carrefour.com.br/busca/lingui%C3%A7a"),
    @("Não encontrado", "https://mercado.carrefour.com.br/busca/lingui%C3%A7a?page=3"),
    @("Ovo Branco Grande Mantiqueira Happy Eggs com 20 Unidades", "https://mercado.carrefour.com.br/ovo-branco-grande-mantiqueira-happy-eggs-com-20-unidades-6403565/p"),
    @("Óleo de Soja Confiare 900ml", "https://mercado.carrefour.com.br/oleo-de-soja-confiare-900ml-3731243/p")
)

$startRow = 178
$runDate = "2025-09-18"
$endRow = $startRow + $data.Length - 1
$dateRange = $ws.Range("D$startRow`:D$endRow")

# Keep the Data column as plain text (matches the existing rows), so Excel
# doesn't auto-convert the "2025-09-18" strings into date serials.
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $runDate
}

# Drop back to the default/Normal style so the new cells don't carry an
# explicit "@" number-format style index (matches the unstyled cells used
# by every other data row in the sheet).
$dateRange.Style = "Normal"
